$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row L1:O1 (reuse existing column labels) + new P1/Q1 headers ---
$ws.Range("L1").Value = "AGO_IND"
$ws.Range("M1").Value = "UNIT_TREATED"
$ws.Range("N1").Value = "AI_RATE_TYPE"
$ws.Range("O1").Value = "SITE_TYPE"
$ws.Range("P1").Value = "MIN_FIXED2"
$ws.Range("Q1").Value = "MAX_FIXED2"

# --- Data rows L2:Q26 (fixed min/max outlier summary, already in final sort order) ---
$ws.Range("L2").Value = "A"
$ws.Range("M2").Value = "A"
$ws.Range("N2").Value = "NORMAL"
$ws.Range("O2").Value = "ALL"
$ws.Range("P2").Value = 0.000000000068443708450083195
$ws.Range("Q2").Value = 3080
$ws.Range("L3").Value = "A"
$ws.Range("M3").Value = "A"
$ws.Range("N3").Value = "MEDIUM"
$ws.Range("O3").Value = "ALL"
$ws.Range("P3").Value = 0.0000000032266319697896301
$ws.Range("Q3").Value = 300
$ws.Range("L4").Value = "A"
$ws.Range("M4").Value = "A"
$ws.Range("N4").Value = "HIGH"
$ws.Range("O4").Value = "ALL"
$ws.Range("P4").Value = 0.0000012935860897065699
$ws.Range("Q4").Value = 3018.1086519114601
$ws.Range("L5").Value = "A"
$ws.Range("M5").Value = "A"
$ws.Range("N5").Value = "ADJUVANT"
$ws.Range("O5").Value = "ALL"
$ws.Range("P5").Value = 0.000000016209751617927999
$ws.Range("Q5").Value = 360
$ws.Range("L6").Value = "A"
$ws.Range("M6").Value = "C"
$ws.Range("N6").Value = "NORMAL"
$ws.Range("O6").Value = "ALL"
$ws.Range("P6").Value = 0.00324675324675324
$ws.Range("Q6").Value = 2360824.7422680398
$ws.Range("L7").Value = "A"
$ws.Range("M7").Value = "C"
$ws.Range("N7").Value = "ADJUVANT"
$ws.Range("O7").Value = "ALL"
$ws.Range("P7").Value = 0.76893939393939303
$ws.Range("Q7").Value = 1238095.23809523
$ws.Range("L8").Value = "A"
$ws.Range("M8").Value = "P"
$ws.Range("N8").Value = "NORMAL"
$ws.Range("O8").Value = "ALL"
$ws.Range("P8").Value = 0.0050000000000000001
$ws.Range("Q8").Value = 40000000
$ws.Range("L9").Value = "A"
$ws.Range("M9").Value = "P"
$ws.Range("N9").Value = "ADJUVANT"
$ws.Range("O9").Value = "ALL"
$ws.Range("P9").Value = 7.6893939393939297
$ws.Range("Q9").Value = 40000000
$ws.Range("L10").Value = "A"
$ws.Range("M10").Value = "U"
$ws.Range("N10").Value = "NORMAL"
$ws.Range("O10").Value = "ALL"
$ws.Range("P10").Value = 0.064935064935064901
$ws.Range("Q10").Value = 41.24
$ws.Range("L11").Value = "A"
$ws.Range("M11").Value = "U"
$ws.Range("N11").Value = "HIGH"
$ws.Range("O11").Value = "ALL"
$ws.Range("P11").Value = 6.6225165562913899
$ws.Range("Q11").Value = 397.67441860465101
$ws.Range("L12").Value = "A"
$ws.Range("M12").Value = "U"
$ws.Range("N12").Value = "ADJUVANT"
$ws.Range("O12").Value = "ALL"
$ws.Range("P12").Value = 15.378787878787801
$ws.Range("Q12").Value = 23.446153846153798
$ws.Range("L13").Value = "N"
$ws.Range("M13").Value = "A"
$ws.Range("N13").Value = "NORMAL"
$ws.Range("O13").Value = "OTHER"
$ws.Range("P13").Value = 0.00000000010600134847166899
$ws.Range("Q13").Value = 708.247422680412
$ws.Range("L14").Value = "N"
$ws.Range("M14").Value = "A"
$ws.Range("N14").Value = "MEDIUM"
$ws.Range("O14").Value = "OTHER"
$ws.Range("P14").Value = 0.000000119092052703144
$ws.Range("Q14").Value = 10000
$ws.Range("L15").Value = "N"
$ws.Range("M15").Value = "A"
$ws.Range("N15").Value = "HIGH"
$ws.Range("O15").Value = "OTHER"
$ws.Range("P15").Value = 0.0000000083213140273522205
$ws.Range("Q15").Value = 10000
$ws.Range("L16").Value = "N"
$ws.Range("M16").Value = "A"
$ws.Range("N16").Value = "ADJUVANT"
$ws.Range("O16").Value = "OTHER"
$ws.Range("P16").Value = 0.000000115037173639728
$ws.Range("Q16").Value = 371.42857142857099
$ws.Range("L17").Value = "N"
$ws.Range("M17").Value = "C"
$ws.Range("N17").Value = "NORMAL"
$ws.Range("O17").Value = "OTHER"
$ws.Range("P17").Value = 0.067044920096464594
$ws.Range("Q17").Value = 810000000
$ws.Range("L18").Value = "N"
$ws.Range("M18").Value = "C"
$ws.Range("N18").Value = "ADJUVANT"
$ws.Range("O18").Value = "OTHER"
$ws.Range("P18").Value = 72.760000000000005
$ws.Range("Q18").Value = 105000000
$ws.Range("L19").Value = "N"
$ws.Range("M19").Value = "P"
$ws.Range("N19").Value = "NORMAL"
$ws.Range("O19").Value = "OTHER"
$ws.Range("P19").Value = 0.0013408984019292901
$ws.Range("Q19").Value = 484150000
$ws.Range("L20").Value = "N"
$ws.Range("M20").Value = "P"
$ws.Range("N20").Value = "MEDIUM"
$ws.Range("O20").Value = "OTHER"
$ws.Range("P20").Value = 10
$ws.Range("Q20").Value = 40000000
$ws.Range("L21").Value = "N"
$ws.Range("M21").Value = "P"
$ws.Range("N21").Value = "HIGH"
$ws.Range("O21").Value = "OTHER"
$ws.Range("P21").Value = 2
$ws.Range("Q21").Value = 200000000
$ws.Range("L22").Value = "N"
$ws.Range("M22").Value = "P"
$ws.Range("N22").Value = "ADJUVANT"
$ws.Range("O22").Value = "OTHER"
$ws.Range("P22").Value = 1.4552
$ws.Range("Q22").Value = 8000000
$ws.Range("L23").Value = "N"
$ws.Range("M23").Value = "U"
$ws.Range("N23").Value = "NORMAL"
$ws.Range("O23").Value = "OTHER"
$ws.Range("P23").Value = 0.670449200964646
$ws.Range("Q23").Value = 10000
$ws.Range("L24").Value = "N"
$ws.Range("M24").Value = "U"
$ws.Range("N24").Value = "MEDIUM"
$ws.Range("O24").Value = "OTHER"
$ws.Range("P24").Value = 48.753894080996801
$ws.Range("Q24").Value = 1000
$ws.Range("L25").Value = "N"
$ws.Range("M25").Value = "U"
$ws.Range("N25").Value = "HIGH"
$ws.Range("O25").Value = "OTHER"
$ws.Range("P25").Value = 41.6666666666666
$ws.Range("Q25").Value = 8000
$ws.Range("L26").Value = "N"
$ws.Range("M26").Value = "U"
$ws.Range("N26").Value = "ADJUVANT"
$ws.Range("O26").Value = "OTHER"
$ws.Range("P26").Value = 727.6
$ws.Range("Q26").Value = 1000

# --- Number formatting: thousands-separated integer format on rate columns + new min/max columns ---
$ws.Range("E1:G36").NumberFormat = "#,##0"
$ws.Range("P1:Q26").NumberFormat = "#,##0"

# --- Column width tweaks (approximate autosize the real Excel session produced) ---
$ws.Columns.Item(5).ColumnWidth = 5.666666666666667
$ws.Columns.Item(6).ColumnWidth = 5.666666666666667
$ws.Columns.Item(7).ColumnWidth = 6.666666666666667
$ws.Columns.Item(14).ColumnWidth = 12.666666666666666
$ws.Columns.Item(16).ColumnWidth = 8.5
$ws.Columns.Item(17).ColumnWidth = 10.333333333333334

# --- Record the 3-key sort applied to build the summary table (data already in final order) ---
$sortRange = $ws.Range("L1:Q26")
$ws.Sort.SortFields.Clear()
[void]$ws.Sort.SortFields.Add($ws.Range("L2:L26"), 0, 1)
[void]$ws.Sort.SortFields.Add($ws.Range("M2:M26"), 0, 1)
[void]$ws.Sort.SortFields.Add($ws.Range("N2:N26"), 0, 2)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Selection / view state ---
[void]$ws.Range("L36").Select()
